$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.200.2.1
$ws.Range("C3").Value = 4665
$ws.Range("D3").Value = 91.59999999999999

# Row 4 - MediaTek Wi-Fi 6 MT7921 Wireless LAN Card - 23.32.2.560
$ws.Range("C4").Value = 346
$ws.Range("D4").Value = 96.09999999999999

# Row 5 - Intel(R) Wi-Fi 6E AX211 160MHz - 23.110.0.5
$ws.Range("C5").Value = 5061
$ws.Range("D5").Value = 97.5

# Row 6 - Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 460

# Row 7 - driver name swapped to AX200 160MHz - 23.60.0.10
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.60.0.10"
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 10

# Row 8 - driver name swapped to AX201 160MHz - 23.60.1.2
$ws.Range("A8").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 167

# Row 9 - driver name swapped to AX201 160MHz - 22.200.2.1
$ws.Range("A9").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.200.2.1"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 50

# Row 10 - Totals
$ws.Range("B10").Value = 126
$ws.Range("C10").Value = 10759

# Row 18 - Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("B18").Value = 449371

# Row 21 - Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8
$ws.Range("B21").Value = 331283

# Row 29 - Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("B29").Value = 77999
